# Commit: "change metadata sheet to isa template"
# Rename the "SwateTemplateMetadata" worksheet to "isa_template".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"
